# Generate Report for Archive
# - Status text "Ready for handoff" -> "In Translation" across all sheets
# - Shrink the now-narrower "Status" columns to match (AutoFit-style resize)

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn / de-de status columns (E and F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($addr in @("E2", "F2", "E3", "F3", "E4", "F4")) {
    $cell = $wsOverview.Range($addr)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}
# Narrow columns E and F now that the status text is shorter.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn / de-de sheets: Status column (C) ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in @("C2", "C3", "C4")) {
        $cell = $ws.Range($addr)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
    # Narrow column C now that the status text is shorter.
    $ws.Columns.Item(3).ColumnWidth = 12.5
}
